$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-22 down to 13-23
$ws.Rows("12:12").Insert()

# Copy the formatting from row 11 (the row above) onto the newly inserted row 12,
# matching Excel's default "insert copies format from above" behaviour
$ws.Range("A11:J11").Copy()
$ws.Range("A12:J12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new project's data (XPROMO)
$ws.Range("A12").Value = "xpromo"
$ws.Range("B12").Value = "Uma rede social de compras coletivas. Explore seu mundo de interesses!"
$ws.Range("C12").Value = 42932
$ws.Range("D12").Value = "DEV"
$ws.Range("E12").Value = "X"
$ws.Range("I12").Value = "MPS et al."
$ws.Range("J12").Value = "X"

# Reset the active cell/selection
$ws.Range("A2").Select()
